# Auto-generated Excel COM-interop edit script
# Applies Spanish (es-us) and English (en-us) localization string updates
# per the commit 'Updating Foreign Language to V74'

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Localization")

$v_C129 = @'
**No ha seleccionado ninguna opción. Vuelva a comenzar y seleccione opciones para cada pregunta, para que yo pueda hacerle recomendaciones.**
'@
$ws.Range("C129").Value = $v_C129

$v_C140 = @'
¡Gracias! El lugar donde se encuentra tiene su propia herramienta de autoevaluación. Por favor [seleccione aquí](
'@
$ws.Range("C140").Value = $v_C140

$v_C144 = @'
**Dígale a un administrador o al personal de enfermería de la escuela o guardería que pudo haber tenido contacto con alguien con COVID-19 presunto.**
'@
$ws.Range("C144").Value = $v_C144

$v_C160 = @'
**Sus síntomas podrían ser causados por el COVID-19.** Si bien la mayoría de las personas se recuperarán y vuelven a su salud normal, algunas tienen síntomas que pueden durar semanas o meses después de tener el COVID-19. Es importante que informe a su proveedor de atención médica sobre sus síntomas porque esto podría afectar sus necesidades médicas en el futuro.
'@
$ws.Range("C160").Value = $v_C160

$v_C168 = @'
**Dígale a un cuidador en su establecimiento que usted pudo haber estado expuesto al COVID-19 en los últimos 14 días.**
'@
$ws.Range("C168").Value = $v_C168

$v_C176 = @'
**Dígale al proveedor de salud ocupacional (o supervisor) de su lugar de trabajo que se está sintiendo enfermo, lo antes posible.**
'@
$ws.Range("C176").Value = $v_C176

$v_C182 = @'
**Aíslese de las demás personas por al menos 5 días.** Para calcular su periodo de aislamiento de 5 días, el día 0 es el primer día de síntomas o el día de la prueba (no el día que recibió el resultado positivo de la prueba) si no tiene síntomas. El día 1 es el primer día entero después de que comenzaron sus síntomas o el primer día entero después del día en que se hizo la prueba. Use una mascarilla de alta calidad si tiene que estar alrededor de otras personas en su casa y vigile sus síntomas.
- **Puede terminar el aislamiento después de 5 días completos si lleva 24 horas sin fiebre sin el uso de medicamentos para reducirla y sus otros síntomas han mejorado.** Si está inmunodeprimido o tiene síntomas de COVID-19 de moderados a graves, deberá aislarse hasta el día 10. Si está inmunodeprimido o se enfermó gravemente, también deberá consultar a un proveedor de atención médica sobre cuándo sería adecuado terminar su aislamiento.
- **Debería seguir usando una mascarilla de alta calidad cuando esté cerca de otras personas en casa y en sitios públicos por otros 5 días (del día 6 al día 10) después de que termine su periodo de aislamiento de 5 días.** Si tiene acceso a pruebas de antígenos, debería considerar usarlas después de salir del aislamiento el día 6. Si tiene dos pruebas negativas secuenciales con una separación de 48 horas (se puede hacer la primera tan pronto como el día 6), puede quitarse la mascarilla antes del día 10. Si los resultados de la prueba de antígenos son positivos, usted todavía podría ser infeccioso y no debería quitarse la mascarilla cuando esté alrededor de otras personas. Siga haciéndose pruebas de antígenos con una separación de al menos 48 horas, hasta que tenga dos resultados negativos secuenciales. Esto podría significar que debe seguir usando una mascarilla y haciéndose pruebas más allá del día 10. No se acerque a las personas que estén inmunodeprimidas o en alto riesgo de enfermarse gravemente, y no vaya a hogares de ancianos y otros entornos de alto riesgo, hasta después de al menos 10 días después de la exposición.
- **No se acerque a las personas que estén inmunodeprimidas o en alto riesgo de enfermarse gravemente, y no vaya a hogares de ancianos y otros entornos de alto riesgo, hasta después de al menos 10 días después de la exposición.**
'@
$ws.Range("C182").Value = $v_C182

$v_C185 = @'
**Los CDC recomiendan que usted se aísle y se mantenga alejado de las demás personas en su establecimiento por al menos 5 días.** Para calcular su periodo de aislamiento de 5 días, el día 0 es el primer día de síntomas o el día de la prueba que dio positivo si no tiene síntomas. El día 1 es el primer día completo después de que presentó síntomas o de la prueba que dio positivo. Si está inmunodeprimido o tiene síntomas de COVID-19 de moderados a graves, deberá aislarse por al menos 10 días y consultar a un proveedor de atención médica sobre cuándo sería adecuado terminar su aislamiento.
Use una mascarilla de alta calidad si tiene que estar alrededor de otras personas en su establecimiento y vigile sus síntomas. Si no puede usar una mascarilla cuando esté alrededor de los demás, debe continuar el aislamiento por 10 días completos. Puede que su establecimiento tenga directrices más específicas. Siga las directrices de su establecimiento cuanto tome precauciones.
'@
$ws.Range("C185").Value = $v_C185

$v_C187 = @'
**Si presenta síntomas nuevos, deberá aislarse y hacerse la prueba de COVID-19. Siga quedándose en casa hasta que tenga los resultados.**
'@
$ws.Range("C187").Value = $v_C187

$v_C188 = @'
**Hasta que pueda recibir los resultados de la prueba de COVID-19, aíslese y manténgase alejado de otras personas por al menos 5 días completos desde que comenzó a tener síntomas.** Si todavía no tiene los resultados de la prueba de COVID-19, use una mascarilla de alta calidad cuando esté alrededor de otras personas en casa y en sitios públicos por otros 5 días. No vaya a sitios donde no pueda usar una mascarilla, evite viajar y estar alrededor de personas que tengan más probabilidad de enfermarse gravemente de COVID-19. Para calcular su periodo de aislamiento de 5 días, el día 0 es el primer día de síntomas. El día 1 es el primer día completo después de que presentó síntomas.
'@
$ws.Range("C188").Value = $v_C188

$v_C190 = @'
**Llame a su proveedor médico, a una línea de asesoría clínica o a un proveedor de telemedicina lo más pronto posible.** También tiene afecciones que podrían ponerlo en un mayor riesgo de enfermarse gravemente.
'@
$ws.Range("C190").Value = $v_C190

$v_C197 = @'
La finalidad del Autoverificador del Coronavirus es ayudarlo a tomar decisiones sobre la obtención de atención médica adecuada. Este sistema no se ha diseñado para fines de diagnóstico o tratamiento de enfermedades, incluido el COVID-19.


Este proyecto ha sido posible mediante la colaboración con la Fundación de los CDC y se ha habilitado a través de la plataforma Azure de Microsoft. La colaboración de los CDC con una organización no federal no implica el respaldo a ningún servicio, producto o empresa específicos.


Para continuar usando esta herramienta, confirme que ha leído y entendido el contenido de este descargo de responsabilidad.


###### ver74 (9.12.2022)
'@
$ws.Range("C197").Value = $v_C197

$v_B198 = @'
**Steps to follow every day:**
- Stay up to date on vaccination, including recommended booster doses. You are up to date if you have completed a primary series and received the most recent booster dose recommended for you by CDC.
- Maintain ventilation improvements.
- Avoid contact with people who have suspected or confirmed COVID-19.
- Follow recommendations for isolation if you have suspected or confirmed COVID-19.
- Follow the recommendations for what to do if you are exposed to someone with COVID-19.
- If you are at high risk of getting very sick, talk with a healthcare provider about additional prevention actions.


**Select the links below for more information on:**
- [COVID-19 symptoms](https://www.cdc.gov/coronavirus/2019-ncov/symptoms-testing/symptoms.html)
- [Post-COVID Conditions](https://www.cdc.gov/coronavirus/2019-ncov/long-term-effects.html)
- [When to get tested](https://www.cdc.gov/coronavirus/2019-ncov/testing/diagnostic-testing.html#who-should-get-tested)
- [Protecting yourself and others from getting sick](https://www.cdc.gov/coronavirus/2019-ncov/prevent-getting-sick/prevention.html)
- [When to isolate and for how long](https://www.cdc.gov/coronavirus/2019-ncov/if-you-are-sick/quarantine.html)
- [What to do if you were exposed to COVID-19](https://www.cdc.gov/coronavirus/2019-ncov/if-you-are-sick/quarantine.html)
- [Taking care of yourself when you are sick](https://www.cdc.gov/coronavirus/2019-ncov/if-you-are-sick/steps-when-sick.html)
- [Taking care of someone else who is sick](https://www.cdc.gov/coronavirus/2019-ncov/if-you-are-sick/care-for-someone.html)
- [Treatments your healthcare provider might recommend if you are sick](https://www.cdc.gov/coronavirus/2019-ncov/your-health/treatments-for-severe-illness.html)
- [Learn about COVID-19 Vaccines](https://www.cdc.gov/coronavirus/2019-ncov/vaccines/index.html)
- [Find COVID-19 vaccine locations near you](https://www.vaccines.gov/)
'@
$ws.Range("B198").Value = $v_B198

$v_C198 = @'
**Medidas para tomar todos los días:**
- Mantenerse al día con las vacunas, incluso las dosis de refuerzo recomendadas. Usted está al día si ha completado una serie primaria y recibido la dosis de refuerzo más reciente que los CDC hayan recomendado para usted.
- Mantener las mejoras en la ventilación.
- Evitar el contacto con las personas con COVID-19 presunto o confirmado.
- Seguir las recomendaciones de aislamiento si usted tiene COVID-19 presunto o confirmado.
- Seguir las recomendaciones acerca de qué hacer si se expone a alguien con COVID-19.
- Si está en alto riesgo de enfermarse gravemente, hable con un proveedor de atención médica acerca de medidas de prevención adicionales.
**Seleccione los enlaces a continuación para obtener más información sobre lo siguiente:**
- [Síntomas del COVID-19](https://espanol.cdc.gov/coronavirus/2019-ncov/symptoms-testing/symptoms.html)
- [Afecciones pos-COVID-19](https://espanol.cdc.gov/coronavirus/2019-ncov/long-term-effects.html)
- [Cuándo hacerse la prueba](https://espanol.cdc.gov/coronavirus/2019-ncov/testing/diagnostic-testing.html)
- [Cómo protegerse a sí mismo y a los demás para que no se enfermen](https://espanol.cdc.gov/coronavirus/2019-ncov/prevent-getting-sick/prevention.html)
- [Cuándo aislarse y por cuánto tiempo](https://espanol.cdc.gov/coronavirus/2019-ncov/if-you-are-sick/quarantine.html)
- [Qué hacer si se expuso al COVID-19](https://espanol.cdc.gov/coronavirus/2019-ncov/if-you-are-sick/quarantine.html)
- [Cómo cuidarse cuando esté enfermo](https://espanol.cdc.gov/coronavirus/2019-ncov/if-you-are-sick/steps-when-sick.html)
- [Cómo cuidar a otra persona que esté enferma](https://espanol.cdc.gov/coronavirus/2019-ncov/if-you-are-sick/care-for-someone.html)
- [Tratamientos que su proveedor de atención médica podría recomendarle si está enfermo](https://espanol.cdc.gov/coronavirus/2019-ncov/your-health/treatments-for-severe-illness.html)
- [Infórmese sobre las vacunas contra el COVID-19](https://espanol.cdc.gov/coronavirus/2019-ncov/vaccines/index.html)
- [Encuentre sitios de vacunación contra el COVID-19 cercanos](https://www.vacunas.gov/)
'@
$ws.Range("C198").Value = $v_C198

$v_B200 = @'
Hi, I’m Clara. I’m going to ask you some questions. I will use your answers to give you advice about the steps you should take to protect yourself and others from COVID-19. If answering for someone else, please respond to all questions as if you are them. If you need to start over, refresh the page and start again.
If you are experiencing a life-threatening emergency, please call 911 immediately.
If you are not experiencing a life-threatening emergency, let’s get started.
'@
$ws.Range("B200").Value = $v_B200

$v_C200 = @'
Hola, me llamo Clara. Le voy a hacer algunas preguntas. Usaré sus respuestas para aconsejarle sobre las medidas que debería tomar para protegerse y proteger a los demás del COVID-19. Si está respondiendo por alguien más, responda todas las preguntas como si fuera esa persona. Si necesita comenzar de nuevo, refresque la página y comience otra vez.
Si tiene una emergencia potencialmente mortal, llame al 911 inmediatamente.
Si no tiene una emergencia potencialmente mortal, comencemos.
'@
$ws.Range("C200").Value = $v_C200

$v_B201 = @'
**CDC recommends these steps to protect you and others from COVID-19:**
- **Stay up to date on vaccination, including recommended booster doses. You are up to date if you have completed a primary series and received the most recent booster dose recommended for you by CDC.**
- **Maintain ventilation improvements.**
- **Avoid contact with people who have suspected or confirmed COVID-19.**
- **Follow recommendations for isolation if you have suspected or confirmed COVID-19.**
- **Follow the recommendations for what to do if you are exposed to someone with COVID-19.**
- **If you are at high risk of getting very sick, talk with a healthcare provider about additional prevention actions.**

'@
$ws.Range("B201").Value = $v_B201

$v_C201 = @'
**Los CDC recomiendan estas medidas para protegerse y proteger a los demás del COVID-19:**
- **Mantenerse al día con las vacunas, incluso las dosis de refuerzo recomendadas. Usted está al día si ha completado una serie primaria y recibido la dosis de refuerzo más reciente que los CDC hayan recomendado para usted.**
- **Mantener las mejoras en la ventilación.**
- **Evitar el contacto con las personas con COVID-19 presunto o confirmado.**
- **Seguir las recomendaciones de aislamiento si usted tiene COVID-19 presunto o confirmado.**
- **Seguir las recomendaciones acerca de qué hacer si se expone a alguien con COVID-19.**
- **Si está en alto riesgo de enfermarse gravemente, hable con un proveedor de atención médica acerca de medidas de prevención adicionales.**
'@
$ws.Range("C201").Value = $v_C201

$v_B202 = @'
Hi, I’m Clara. I’m going to ask you some questions. I will use your answers to give you advice about the steps you should take to protect yourself and others from COVID-19. If answering for someone else, please respond to all questions as if you are them. If you need to start over, refresh the page and start again.
If you are experiencing a life-threatening emergency, please call 911 immediately.
If you are not experiencing a life-threatening emergency, let’s get started.
**CDC recommends these steps to protect you and others from COVID-19:**
- **Stay up to date on vaccination, including recommended booster doses. You are up to date if you have completed a primary series and received the most recent booster dose recommended for you by CDC.**
- **Maintain ventilation improvements.**
- **Avoid contact with people who have suspected or confirmed COVID-19.**
- **Follow recommendations for isolation if you have suspected or confirmed COVID-19.**
- **Follow the recommendations for what to do if you are exposed to someone with COVID-19.**
- **If you are at high risk of getting very sick, talk with a healthcare provider about additional prevention actions.**
'@
$ws.Range("B202").Value = $v_B202

$v_C202 = @'
Hola, me llamo Clara. Le voy a hacer algunas preguntas. Usaré sus respuestas para aconsejarle sobre las medidas que debería tomar para protegerse y proteger a los demás del COVID-19. Si está respondiendo por alguien más, responda todas las preguntas como si fuera esa persona. Si necesita comenzar de nuevo, refresque la página y comience otra vez.
Si tiene una emergencia potencialmente mortal, llame al 911 inmediatamente.
Si no tiene una emergencia potencialmente mortal, comencemos.
**Los CDC recomiendan estas medidas para protegerse y proteger a los demás del COVID-19:**
- **Mantenerse al día con las vacunas, incluso las dosis de refuerzo recomendadas. Usted está al día si ha completado una serie primaria y recibido la dosis de refuerzo más reciente que los CDC hayan recomendado para usted.**
- **Mantener las mejoras en la ventilación.**
- **Evitar el contacto con las personas con COVID-19 presunto o confirmado.**
- **Seguir las recomendaciones de aislamiento si usted tiene COVID-19 presunto o confirmado.**
- **Seguir las recomendaciones acerca de qué hacer si se expone a alguien con COVID-19.**
- **Si está en alto riesgo de enfermarse gravemente, hable con un proveedor de atención médica acerca de medidas de prevención adicionales.**
'@
$ws.Range("C202").Value = $v_C202

$v_C208 = @'
**Debido a que algunos de los síntomas de la influenza y del COVID-19 son similares, es recomendable que le pregunte a su proveedor si le aconseja la prueba o tratamiento para la influenza.** Consulte esta [página web de los CDC](https://espanol.cdc.gov/flu/symptoms/flu-vs-covid19.htm) para obtener más información sobre el COVID-19 y la influenza.
'@
$ws.Range("C208").Value = $v_C208

$v_C215 = @'
**Independientemente de si se ha vacunado o si tuvo una infección anterior, usted podría necesitar hacerse la prueba de nuevo si no hay otra causa identificada para sus síntomas.** Considere volver a hacerse la prueba de la infección por SARS-CoV-2 debido a la posibilidad de reinfección.
'@
$ws.Range("C215").Value = $v_C215

$v_C217 = @'
**Si tiene preguntas, hable con su proveedor de atención médica acerca de los resultados de su prueba y el tipo de prueba que se hizo para entender lo que significan esos resultados.**
'@
$ws.Range("C217").Value = $v_C217

$v_C218 = @'
**Si no se ha identificado otra causa para sus síntomas, considere volver a hacerse la prueba.** Si tiene preguntas, hable con su proveedor de atención médica acerca de los resultados de su prueba y el tipo de prueba que se hizo para entender lo que significan esos resultados.
'@
$ws.Range("C218").Value = $v_C218

$v_C224 = @'
Si sigue sin síntomas de COVID-19, no necesita ponerse en cuarentena ni hacerse la prueba a menos que se lo recomiende o exija su proveedor de atención médica, empleador o funcionario de salud pública.
'@
$ws.Range("C224").Value = $v_C224

$v_C225 = @'
Según las respuestas que nos ha dado, no necesita ponerse en cuarentena ni hacerse la prueba a menos que se lo recomiende o exija su proveedor de atención médica, empleador o funcionario de salud pública.
'@
$ws.Range("C225").Value = $v_C225

$v_C231 = @'
**Debido a que usted tiene síntomas de COVID-19, deberá aislarse de los demás inmediatamente.** Deberá seguir aislándose y usar una mascarilla de alta calidad si tiene que estar alrededor de otras personas, hasta que pueda recibir los resultados de la prueba de COVID-19.
'@
$ws.Range("C231").Value = $v_C231

$v_C232 = @'
**Debido a que usted estuvo o puede haber estado expuesto a alguien con COVID-19, los CDC recomiendan que use una mascarilla de alta calidad cuando esté alrededor de otras personas por 10 días completos después de su exposición.  Debería hacerse la prueba después de 5 días completos de haber estado expuesto a alguien con COVID-19, aunque no presente síntomas. Si presenta síntomas, aíslese inmediatamente y hágase la prueba.** Por 10 días completos después de la exposición, también deberá tomar precauciones adicionales con los viajes o si va a estar alrededor de personas que tengan más probabilidad de enfermarse gravemente de COVID-19.
'@
$ws.Range("C232").Value = $v_C232

$v_C249 = @'
["Menores de 2 años","2-4 años","5-9 años","10-12 años","13-17 años","18-29 años","30-39 años","40-49 años","50-59 años","60-64 años","65-69 años","70-79 años","80 años o más"]
'@
$ws.Range("C249").Value = $v_C249

$v_C257 = @'
Asma
'@
$ws.Range("C257").Value = $v_C257

$v_C259 = @'
Fibrosis quística
'@
$ws.Range("C259").Value = $v_C259

$v_C284 = @'
Estado de inmunodepresión (sistema inmunitario debilitado) o uso de medicamentos inmunodepresores
'@
$ws.Range("C284").Value = $v_C284

$v_C356 = @'
¿Tiene alguno de estos síntomas potencialmente mortales?
- Dificultad para respirar
- Dolor o presión persistentes en el pecho
- Confusión de aparición reciente
- No poder despertarse o permanecer despierto
- Color pálido, gris o azulado de la piel, los labios o el lecho de las uñas, dependiendo del tono de piel

*Esta lista no tiene todos los síntomas posibles. Llame a un proveedor de atención médica si tiene algún síntoma grave o que le preocupe.
'@
$ws.Range("C356").Value = $v_C356

$v_C375 = @'
En las últimas dos semanas, ¿ha estado [expuesto al](https://espanol.cdc.gov/coronavirus/2019-ncov/your-health/risks-exposure.html) COVID-19?
'@
$ws.Range("C375").Value = $v_C375

$v_C381 = @'
En las últimas dos semanas, ¿ha trabajado o ha servido como voluntario en un establecimiento de atención médica? Los establecimientos de atención médica incluyen hospitales, centros médicos o dentales, establecimientos de cuidados a largo plazo u hogares de ancianos.
'@
$ws.Range("C381").Value = $v_C381

$v_C385 = @'
¿Vive en un establecimiento de cuidados a largo plazo, hogar de ancianos, centro correccional o refugio para personas sin hogar?
'@
$ws.Range("C385").Value = $v_C385

$v_C397 = @'
**¿Está al día con sus vacunas contra el COVID-19 [es decir que ha completado una serie primaria y recibido la dosis de refuerzo más reciente que los CDC recomendaron para usted?](https://espanol.cdc.gov/coronavirus/2019-ncov/vaccines/stay-up-to-date.html)**
'@
$ws.Range("C397").Value = $v_C397

